# Applies scheduled-runner price/profit refresh to the Leve profit sheets.
# Data comes from the diff: columns H-N (currentAveragePrice*, LevePrice*, LeveProfit*)
# are refreshed per-row; a couple of rows gain/lose a trailing cell where the
# upstream feed newly has/lacks a HQ or NQ price.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1357.2632
$ws.Range("J17").Value = 1377.5272
$ws.Range("L17").Value = 4132.5816
$ws.Range("N17").Value = -4468.5816
$ws.Range("H46").Value = 4698.75
$ws.Range("I46").Value = 4159.6665
$ws.Range("K46").Value = 12478.9995
$ws.Range("M46").Value = -12359.9995
$ws.Range("H60").Value = 4698.75
$ws.Range("I60").Value = 4159.6665
$ws.Range("K60").Value = 12478.9995
$ws.Range("M60").Value = -11994.9995
$ws.Range("H62").Value = 90965910
$ws.Range("I62").Value = 200018350
$ws.Range("K62").Value = 200018350
$ws.Range("M62").Value = -200017726
$ws.Range("H65").Value = 90965910
$ws.Range("I65").Value = 200018350
$ws.Range("K65").Value = 1000091750
$ws.Range("M65").Value = -1000088630
$ws.Range("H98").Value = 4349
$ws.Range("I98").Value = 4606.154
$ws.Range("K98").Value = 4606.154
$ws.Range("M98").Value = -3108.154
$ws.Range("H103").Value = 1400.8462
$ws.Range("I103").Value = 614.5
$ws.Range("J103").Value = 1543.8182
$ws.Range("K103").Value = 1843.5
$ws.Range("L103").Value = 4631.4546
$ws.Range("M103").Value = -1257.5
$ws.Range("N103").Value = -5803.4546
$ws.Range("H112").Value = 5305.814
$ws.Range("I112").Value = 887.6
$ws.Range("J112").Value = 5887.1577
$ws.Range("K112").Value = 2662.8
$ws.Range("L112").Value = 17661.4731
$ws.Range("M112").Value = -1554.8
$ws.Range("N112").Value = -19877.4731
$ws.Range("H122").Value = 4349
$ws.Range("I122").Value = 4606.154
$ws.Range("K122").Value = 13818.462
$ws.Range("M122").Value = -11368.462

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 215
$ws.Range("I4").Value = 174.28572
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 174.28572
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = -58.28572
$ws.Range("N4").Value = -732
$ws.Range("H43").Value = 50000
$ws.Range("J43").Value = 50000
$ws.Range("L43").Value = 50000
$ws.Range("N43").Value = -50626
$ws.Range("H46").Value = 5499.5
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 5499.5
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 5499.5
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -6137.5
$ws.Range("H132").Value = 5989.5713
$ws.Range("J132").Value = 10618.182
$ws.Range("L132").Value = 31854.546
$ws.Range("N132").Value = -36914.546
$ws.Range("H133").Value = 94629.7
$ws.Range("J133").Value = 94629.7
$ws.Range("L133").Value = 94629.7
$ws.Range("N133").Value = -99689.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4667.228
$ws.Range("J134").Value = 8730.708000000001
$ws.Range("L134").Value = 26192.124
$ws.Range("N134").Value = -31262.124

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 147.5
$ws.Range("I7").Value = 33.166668
$ws.Range("J7").Value = 233.25
$ws.Range("K7").Value = 33.166668
$ws.Range("L7").Value = 233.25
$ws.Range("M7").Value = 79.833332
$ws.Range("N7").Value = -459.25
$ws.Range("H31").Value = 7271.7383
$ws.Range("I31").Value = 2692.5881
$ws.Range("K31").Value = 2692.5881
$ws.Range("M31").Value = -2397.5881
$ws.Range("H34").Value = 7271.7383
$ws.Range("I34").Value = 2692.5881
$ws.Range("K34").Value = 2692.5881
$ws.Range("M34").Value = -2490.5881

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 77036.89999999999
$ws.Range("I2").Value = 13679.363
$ws.Range("J2").Value = 251270.12
$ws.Range("K2").Value = 82076.178
$ws.Range("L2").Value = 1507620.72
$ws.Range("M2").Value = -81963.178
$ws.Range("N2").Value = -1507846.72
$ws.Range("H92").Value = 5918437
$ws.Range("J92").Value = 6994352
$ws.Range("L92").Value = 20983056
$ws.Range("N92").Value = -20985552
$ws.Range("H129").Value = 11179362
$ws.Range("I129").Value = 498
$ws.Range("J129").Value = 13974077
$ws.Range("K129").Value = 1494
$ws.Range("L129").Value = 41922231
$ws.Range("M129").Value = 3506
$ws.Range("N129").Value = -41932231

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4432.6113
$ws.Range("I126").Value = 4412.467
$ws.Range("J126").Value = 4533.3335
$ws.Range("K126").Value = 13237.401
$ws.Range("L126").Value = 13600.0005
$ws.Range("M126").Value = -10767.401
$ws.Range("N126").Value = -18540.0005
$ws.Range("H141").Value = 50999.75
$ws.Range("J141").Value = 50999.75
$ws.Range("L141").Value = 50999.75
$ws.Range("N141").Value = -61359.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2843.1724
$ws.Range("I46").Value = 501
$ws.Range("J46").Value = 3735.4285
$ws.Range("K46").Value = 501
$ws.Range("L46").Value = 3735.4285
$ws.Range("M46").Value = -313
$ws.Range("N46").Value = -4111.4285
$ws.Range("H82").Value = 1490
$ws.Range("I82").Value = 949
$ws.Range("J82").Value = 2572
$ws.Range("K82").Value = 949
$ws.Range("L82").Value = 2572
$ws.Range("M82").Value = -588
$ws.Range("N82").Value = -3294
$ws.Range("H85").Value = 1490
$ws.Range("I85").Value = 949
$ws.Range("J85").Value = 2572
$ws.Range("K85").Value = 949
$ws.Range("L85").Value = 2572
$ws.Range("M85").Value = 299
$ws.Range("N85").Value = -5068
$ws.Range("H122").Value = 8000.5
$ws.Range("I122").Value = 4999.5
$ws.Range("K122").Value = 14998.5
$ws.Range("M122").Value = -12548.5
$ws.Range("H132").Value = 10210474
$ws.Range("I132").Value = 22729992
$ws.Range("K132").Value = 68189976
$ws.Range("M132").Value = -68187446

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H102").Value = 53139.5
$ws.Range("J102").Value = 53139.5
$ws.Range("L102").Value = 53139.5
$ws.Range("N102").Value = -59629.5
$ws.Range("H106").Value = 40377
$ws.Range("J106").Value = 40377
$ws.Range("L106").Value = 40377
$ws.Range("N106").Value = -42901
$ws.Range("H119").Value = 56958
$ws.Range("J119").Value = 56958
$ws.Range("L119").Value = 56958
$ws.Range("N119").Value = -66634
$ws.Range("H122").Value = 169350.08
$ws.Range("I122").Value = 268476.28
$ws.Range("J122").Value = 4139.778
$ws.Range("K122").Value = 805428.8400000001
$ws.Range("L122").Value = 12419.334
$ws.Range("M122").Value = -802978.8400000001
$ws.Range("N122").Value = -17319.334
